# Updated symbol list (coin price / 1h volume change) for the 14-1-2023
# 10:00 snapshot, as published by the GitHub Actions scraper job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (Price in column D, Volume(1h) in column E).
$updates = @{
    'D2'  = '302.96'
    'E2'  = '5.70%'
    'D3'  = '31.90'
    'E3'  = '8.08%'
    'D4'  = '5.230'
    'E4'  = '2.47%'
    'D5'  = '0.07262'
    'E5'  = '8.13%'
    'D6'  = '7.787'
    'E6'  = '6.07%'
    'D7'  = '3.750'
    'E7'  = '8.87%'
    'D8'  = '1.460'
    'E8'  = '5.44%'
    'D9'  = '0.9104'
    'E9'  = '-1.05%'
    'D10' = '0.01662'
    'E10' = '2,479.46%'
    'D11' = '0.1679'
    'E11' = '5.09%'
    'D12' = '0.07409'
    'E12' = '8.91%'
    'D13' = '0.07986'
    'E13' = '3.68%'
    'D14' = '0.02982'
    'E14' = '1.82%'
    'D15' = '0.09930'
    'E15' = '10.61%'
    'D16' = '0.001509'
    'E16' = '-4.98%'
    'D17' = '0.04562'
    'E17' = '2.28%'
    'D18' = '0.006558'
    'E18' = '4.77%'
    'D19' = '3.491'
    'E19' = '1.07%'
    'D20' = '2.226'
    'E20' = '-0.09%'
    'D21' = '0.3330'
    'E21' = '3.90%'
    'D22' = '0.1322'
    'E22' = '0.90%'
    'D23' = '4.293'
    'E23' = '5.48%'
    'D25' = '0.001229'
    'E25' = '2.86%'
    'D26' = '0.004414'
    'E26' = '6.88%'
    'D27' = '0.0001314'
    'E27' = '9.61%'
    'E28' = '8.50%'
    'D40' = '0.04481'
    'E40' = '4.57%'
    'D41' = '0.006972'
    'E41' = '3.41%'
    'D42' = '0.1342'
    'E42' = '8.22%'
    'D43' = '0.002416'
    'E43' = '7.98%'
    'D44' = '0.01274'
    'E44' = '6.53%'
    'D45' = '0.00006094'
    'E45' = '7.10%'
    'D47' = '0.01307'
    'E47' = '-13.01%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force the literal text (Excel would otherwise parse these
    # numeric/percent-looking strings as numbers), then drop the
    # resulting "Text" number-format style so the cell stays on the
    # sheet's default (unstyled) format, matching the source data.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.ClearFormats()
}
